# Rewrites Review_295.docx from the "LLMs Will Always Hallucinate" review
# into the "Learning to reason with LLMs" (o1) review:
#  - paragraphs 1-10 (1-based) get their text replaced in place (this also
#    strips the leading <w:br/> that paragraph 3 used to carry, since we
#    overwrite the whole paragraph range rather than just the run text)
#  - the ten paragraphs that used to enumerate hallucination details are
#    deleted outright
#  - the final paragraph's arxiv link is swapped for the openai link
$d = $word.ActiveDocument

$t0 = @"
⚡️🚀המאמר היומי של מייק 12.09.24: ⚡️🚀
"@
$p = $d.Paragraphs.Item(1)
$r = $p.Range
$r.MoveEnd(1, -1) | Out-Null
$r.Text = $t0

$t1 = @"
Learning to reason with LLMs
"@
$p = $d.Paragraphs.Item(2)
$r = $p.Range
$r.MoveEnd(1, -1) | Out-Null
$r.Text = $t1

$t2 = @"
היום במקום הסקירה אשתף איתכם את מחשבותיי על המודל החדש של openai שקיבל שם o1. אני בדרך כלל נמנע מלהגיב ולכתוב פוסטים על כל מודל חדש שמנצח את כל ה-benchmarks בעולם אבל הפעם אחרוג ממנהגי. ולא מהסיבה שמהמודל הזה השאיר אבק לרוב ה-benchmarks אלא בגלל שאני זיהיתי כאן שינוי מסוים בפרדיגמה בעולם ה-llms.
"@
$p = $d.Paragraphs.Item(3)
$r = $p.Range
$r.MoveEnd(1, -1) | Out-Null
$r.Text = $t2

$t3 = @"
השינוי בפרדיגמה בא בדמות של שינוי היחס בכמות הקומפיט המוקדש ללמידה ולהסקה (אינפרנס). אנחנו רגילים למודל שמצריכים כמות אדירה של קומפיוט במהלך הלמידה (אימון מקדים, SFT, יישור המודל וכדומה) כאשר האינפרנס הוא די זול (כמובן יחסית לאימון כי גם בהסקה יש עלויות די גבוהות בשל עצמם). O1 לעומות זאת מאתגר את ההנחה הזו ושואל את השאלה: האם זה אופטימלי? אולי אנו צריכים לאמן את המודל שלנו פחות ולהשקיע יותר קומפיט בהסקה.
"@
$p = $d.Paragraphs.Item(4)
$r = $p.Range
$r.MoveEnd(1, -1) | Out-Null
$r.Text = $t3

$t4 = @"
לפני כמה זמן סקרתי מאמר שדי שינה (או לכל הפחות רענן) את תפיסתי בעניין זה (Scaling LLM Test-Time Compute Optimally can be More Effective than Scaling Model Parameters). המאמר הזה היה של deepmind אולם הייתה לי תחושה שהם לא היחידים שהגיעו לתובנה הדי לא טריוויאלית הזה. 
"@
$p = $d.Paragraphs.Item(5)
$r = $p.Range
$r.MoveEnd(1, -1) | Out-Null
$r.Text = $t4

$t5 = @"
בעקרון הכל מסתכם לשתי הנקודות הבאות:
"@
$p = $d.Paragraphs.Item(6)
$r = $p.Range
$r.MoveEnd(1, -1) | Out-Null
$r.Text = $t5

$t6 = @"
אולי אתה לא צריך מודל שפה ענק להסקה. חלק ניכר מהפרמטרים כנראה ממשמשים לאחסון עובדות, כדי שהמודל לא ידבר שטויות לשאלות לידע כללי (כמו מתי נולד מוצרט). לדעתי ניתן להפריד בין הסקה לידע, כלומר אפשר להסתפק ב"ליבה להסקה" קטנה שיודעת איך להשתמש בכלים כמו וולפרם, בראוזר ובודק קוד כלומר המשימות הדורשות סוג של ידע עובדתי (ידע בשפת תכנות). ככה ניתן להפחית את כמות החישוב המוקדשת לאימון המוקדם.
"@
$p = $d.Paragraphs.Item(7)
$r = $p.Range
$r.MoveEnd(1, -1) | Out-Null
$r.Text = $t6

$t7 = @"
כמות משמעותית של קומפיט מועברת להסקה בזמן הרצת המודל ולא לאימון המודל. ניתן לחשוב על מודלי שפה בתור סימולטורים מבוססי טקסט. על ידי הרצת תרחישים ואסטרטגיות רבות (גנרוט טקסט), המודל יגיע בסופו של דבר לפתרונות reasoning טובים. התהליך בחירת הפתרון נראה די דומה לבעיות שנחקרו היטב כמו חיפוש העץ של מונטה קרלו  (MCTS) ב-AlphaGo. 
"@
$p = $d.Paragraphs.Item(8)
$r = $p.Range
$r.MoveEnd(1, -1) | Out-Null
$r.Text = $t7

$t8 = @"
כמובן שאם יש שימוש בטכניקות כמו MCTS אנו צריכים את פונקציית ה-reward. בניית פונקצייה כזו היא לא טריוויאלית כאן כי אין לנו דרך טובה (אלא אם כן יש לנו דאטהסט reasoning מגוון ועצום שניתן לאמן עליו מודל כזה) לשערך את איכות ה-reasoning. כמובן שניתן לנצל מודלי שפה אחרים, בדיקות עצמיות על ידי מודלי שפה וכדומה אבל עדיין לא ברור ב-100% איך לעשות את זה (ד״א אני בכלל לא בטוח שהם השתמשו ב-mcts). אולי הם פיתחו שיטה מגניבה לעקוף את ה-reward כמו שנעשה ב-dpo וב-orpo שעשו זאת עבור ppo -אין לדעת. 
"@
$p = $d.Paragraphs.Item(9)
$r = $p.Range
$r.MoveEnd(1, -1) | Out-Null
$r.Text = $t8

$t9 = @"
בקיצור מחכה לדוח הטכני שבתקווה ישפוך אור על הסיפור הזה (גם בזה אני לא בטוח בכלל)....
"@
$p = $d.Paragraphs.Item(10)
$r = $p.Range
$r.MoveEnd(1, -1) | Out-Null
$r.Text = $t9

# Delete the large block of paragraphs that is removed entirely
$delStart = $d.Paragraphs.Item(11).Range.Start
$delEnd = $d.Paragraphs.Item(20).Range.End
$delRange = $d.Range($delStart, $delEnd)
$delRange.Delete() | Out-Null

$tlink = @"
https://openai.com/index/learning-to-reason-with-llms/
"@
$p = $d.Paragraphs.Item(11)
$r = $p.Range
$r.MoveEnd(1, -1) | Out-Null
$r.Text = $tlink

Write-Output $d.Paragraphs.Count
